# Add the new sheet "2020年3月“个人成长基金”项目凭证详情" at the end of the
# workbook, mirroring the layout/styling of the existing "2020年3月凭证列表"
# sheet (same template workbook authored the style table, so we reuse its
# cell formats via copy/paste-format rather than re-declaring fonts/borders).

$wb = $excel.ActiveWorkbook

# Template sheet whose cell formats (style ids 2/3/4) we will clone.
$tpl = $wb.Worksheets.Item("2020年3月凭证列表")

# New sheet goes after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "2020年3月“个人成长基金”项目凭证详情"

# --- clone cell formatting from the template sheet -------------------------
# Title row (style id 2: big font, centered, wrapped, no border)
$tpl.Range("A1:K1").Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)   # xlPasteFormats

# Header / text cells (style id 3: bordered, centered, wrapped)
$tpl.Range("B2:G2").Copy()
$ws.Range("B2:G2").PasteSpecial(-4122)

$tpl.Range("B3:C3").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)

$tpl.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

$tpl.Range("B3").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Numeric / bordered-only cells (style id 4)
$tpl.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$tpl.Range("H3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$tpl.Range("H3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$tpl.Range("H3").Copy()
$ws.Range("G3").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- column widths -----------------------------------------------------
$ws.Columns("B").ColumnWidth = 25
$ws.Columns("E").ColumnWidth = 40

# --- values --------------------------------------------------------------
$ws.Range("A1").Value = "2020年3月“个人成长基金”项目凭证详情"

$ws.Range("B2").Value = "摘要"
$ws.Range("C2").Value = "类型"
$ws.Range("D2").Value = "金额"
$ws.Range("E2").Value = "凭证号"
$ws.Range("F2").Value = "小计"
$ws.Range("G2").Value = "合计"

$ws.Range("A3").Value = "支出"
$ws.Range("B3").Value = "个人使用"
$ws.Range("C3").Value = "吃饭"
$ws.Range("D3").Value = 50
$ws.Range("E3").Value = "ABCD"
$ws.Range("F3").Value = 50
$ws.Range("G3").Value = 50

# --- merges ---------------------------------------------------------------
$ws.Range("A1:K1").Merge()
$ws.Range("A3:A3").Merge()

# Row 1 needs a big-font row but no explicit stored row height (matches the
# rest of the workbook, which relies on the sheet-wide default).
$ws.Rows("1:1").AutoFit()
